$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "BQ - Bloc Québécois (Bloc Québécois, BQ)"
$ws.Range("C1").Value = "Ind - Independents (Independents, Ind)"
$ws.Range("D1").Value = "Lib - Liberal  (Liberal , Lib)"
$ws.Range("E1").Value = "ND - New Democratic (New Democratic, ND)"
$ws.Range("F1").Value = "PC - Progressive Conservative  (Progressive Conservative , PC)"
$ws.Range("G1").Value = "RPC/RP - Reform Party of Canada / Canadian Alliance (Reform Party of Canada / Canadian Alliance, RPC/RP)"
$ws.Range("H1").Value = "None - No-Affiliation (No-Affiliation, None)"
$ws.Range("I1").Value = "Con - Conservative (Conservative, Con)"
$ws.Range("J1").Value = "GP - Green Party of Canada (Green Party of Canada, GP)"
